# Marble racing season 2 scoring sheet: "Added autotype for marbles on
# stream scoring system" - fill in the previously-pending Race 1 result
# for MoscaMye and correct the Race 3 / Race 4 results + recompute the
# Total Points for both racers, on all three (differently-sorted) sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Sorted by points", "Sorted by username", "Sorted by flair")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Figure out which row holds which racer on this sheet (rows are
    # sorted differently per tab) by checking column B ("Username").
    if ($ws.Range("B2").Value() -eq "Zokalyx") {
        $zokRow = 2
        $mosRow = 3
    } else {
        $zokRow = 3
        $mosRow = 2
    }

    # --- Zokalyx row -----------------------------------------------
    $ws.Cells.Item($zokRow, 4).Value = 184          # D: Total Points
    $ws.Cells.Item($zokRow, 7).Value = "2nd (+25)"  # G: Race 3
    $ws.Cells.Item($zokRow, 8).Value = "2nd (+99)"  # H: Race 4

    # --- MoscaMye row ------------------------------------------------
    $ws.Cells.Item($mosRow, 4).Value = 180           # D: Total Points

    # E: Race 1 - was blank/pending (grey-highlighted cell); now typed in.
    $eCell = $ws.Cells.Item($mosRow, 5)
    $eCell.Value = "2nd (+25)"
    $eCell.Style = "Normal"
    $eCell.HorizontalAlignment = -4108               # xlCenter - match the rest of the row

    $ws.Cells.Item($mosRow, 7).Value = "1st (+30)"   # G: Race 3
    $ws.Cells.Item($mosRow, 8).Value = "1st (+100)"  # H: Race 4
}
